$wb = $excel.ActiveWorkbook

# Sheet 1 (index 1): 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 265
$ws1.Range("F7").Value = 13110
$ws1.Range("F10").Value = 275
$ws1.Range("F11").Value = 3630
$ws1.Range("F13").Value = 6621
$ws1.Range("F16").Value = 3488
$ws1.Range("F27").Value = 3347
$ws1.Range("F29").Value = 1911
$ws1.Range("F30").Value = 106
$ws1.Range("F32").Value = 6802
$ws1.Range("F34").Value = 1492
$ws1.Range("F35").Value = 2019
$ws1.Range("F38").Value = 1064
$ws1.Range("F40").Value = 219
$ws1.Range("F43").Value = 1147
$ws1.Range("F46").Value = 1810
$ws1.Range("F47").Value = 67
$ws1.Range("F49").Value = 1180

# Sheet 3 (index 3): 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 455
$ws3.Range("F3").Value = 624

# Sheet 4 (index 4): 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 455
$ws4.Range("F7").Value = 624
$ws4.Range("F8").Value = 265
$ws4.Range("F9").Value = 13110
$ws4.Range("F13").Value = 275
$ws4.Range("F14").Value = 3630
$ws4.Range("F16").Value = 3488
$ws4.Range("F27").Value = 3348
$ws4.Range("F28").Value = 3348
$ws4.Range("F30").Value = 1911
$ws4.Range("F31").Value = 106
$ws4.Range("F33").Value = 6802
$ws4.Range("F36").Value = 1493
$ws4.Range("F37").Value = 2019
$ws4.Range("F41").Value = 1064
$ws4.Range("F42").Value = 219
$ws4.Range("F47").Value = 1810
$ws4.Range("F48").Value = 67

Write-Host "Done applying updates."
